$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Paul_example" worksheet as the last (4th) sheet.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Paul_example"

# ---------------------------------------------------------------------------
# 2. Populate the occurrence-matrix data (header + 11 data rows).
# ---------------------------------------------------------------------------
$headers = @("Row Labels","pH4","pH4.5","pH5","pH5.5","pH6","pH6.5","pH7")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$data = @(
    @("F.Otu182", 1, 1, 0, 0, 0, 0, 0, "k__Fungi;p__Ascomycota;c__Leotiomycetes;o__Helotiales;f__Helotiales_fam_Incertae_sedis;g__Leptodontidium;s__unidentified"),
    @("F.Otu157", 0, 0, 0, 0, 0, 1, 1, "k__Fungi;p__Ascomycota;c__Leotiomycetes;o__unidentified;f__unidentified;g__unidentified;s__unidentified"),
    @("F.Otu178", 0, 0, 0, 0, 1, 1, 0, "k__Fungi;p__Ascomycota;c__Sordariomycetes;o__Hypocreales;f__Hypocreaceae;g__Trichoderma;s__Trichoderma_paraviridescens"),
    @("F.Otu201", 0, 0, 0, 1, 1, 0, 0, "k__Fungi;p__Ascomycota;c__Sordariomycetes;o__Sordariales;f__unidentified;g__unidentified;s__unidentified"),
    @("F.Otu128", 0, 0, 1, 0, 0, 0, 1, "k__Fungi;p__Basidiomycota;c__Microbotryomycetes;o__Leucosporidiales;f__Leucosporidiaceae;g__Leucosporidium;s__Leucosporidium_fragarium"),
    @("F.Otu183", 0, 0, 0, 0, 1, 0, 1, "k__Fungi;p__Basidiomycota;c__Tremellomycetes;o__Tremellales;f__Rhynchogastremataceae;g__Papiliotrema;s__Papiliotrema_frias"),
    @("F.Otu238", 0, 0, 0, 0, 1, 1, 0, "k__unidentified;p__unidentified;c__unidentified;o__unidentified;f__unidentified;g__unidentified;s__unidentified"),
    @("B.Otu4918", 0, 1, 1, 0, 0, 0, 0, "p__Chloroflexi; c__Anaerolineae; o__SBR1031; f__uncultured bacterium; g__; s__"),
    @("B.Otu636", 0, 0, 0, 0, 1, 1, 0, "p__Chloroflexi; c__Ktedonobacteria; o__Ktedonobacterales; f__Ktedonobacteraceae; g__JG30a-KF-32; s__uncultured Ktedobacteria bacterium"),
    @("B.Otu368", 0, 0, 1, 1, 0, 0, 0, "p__Firmicutes; c__Bacilli; o__Bacillales; f__Bacillaceae; g__Bacillus; s__Bacillus oceanisediminis"),
    @("B.Otu1321", 0, 0, 1, 0, 0, 1, 0, "p__Planctomycetota; c__Planctomycetes; o__Gemmatales; f__Gemmataceae; g__uncultured; s__uncultured bacterium")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

$ws.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------------------
# 2b. Record the sort that produced this row order (by taxonomy, column I).
# ---------------------------------------------------------------------------
$sortRange = $ws.Range("A2:I22")
$keyRange = $ws.Range("I2:I22")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 3. Conditional formatting (colour scale) on B2:H12.
# ---------------------------------------------------------------------------
$cfRange = $ws.Range("B2:H12")
$cf = $cfRange.FormatConditions.AddColorScale(2)
$cf.ColorScaleCriteria.Item(1).Type = 1
$cf.ColorScaleCriteria.Item(1).FormatColor.Color = 16776444
$cf.ColorScaleCriteria.Item(2).Type = 2
$cf.ColorScaleCriteria.Item(2).FormatColor.Color = 8109667

# ---------------------------------------------------------------------------
# 4. View settings for the new sheet.
# ---------------------------------------------------------------------------
$ws.Range("K17").Select()

# ---------------------------------------------------------------------------
# 5. Adjust the other sheets.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("several_methods")
$ws1.Range("C23").Select()

$ws2 = $wb.Worksheets.Item("several_groups")
$ws2.Range("A1:D29").Select()

$ws.Select()

$wb.Windows.Item(1).WindowState = -4143
